$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.446.47"
$ws.Range("E2").Value = "  -1.10%  "

$ws.Range("D3").Value = "2.524.62"
$ws.Range("E3").Value = "  -0.23%  "

$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.80%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.14%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.573"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.74%  "

$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.527"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.45%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.79"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.35%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0806"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.19%  "

$ws.Range("E12").Value = "  -0.15%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.52"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.09%  "

$ws.Range("D14").Value = "2.912.70"
$ws.Range("E14").Value = "  -0.39%  "

$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.39"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.41%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.505.50"
$ws.Range("E16").Value = "  +0.12%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.844"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.50%  "

$ws.Range("D18").Value = "42.476.43"
$ws.Range("E18").Value = "  -1.05%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.99"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.31%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.84%  "

$ws.Range("D21").Value = "0.0₃0958"
$ws.Range("E21").Value = "  -2.74%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.98"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.30%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "250.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.89%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.87%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.40%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.98%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.996"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.51%  "

$ws.Range("E28").Value = "  +2.47%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.75%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.90"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.46%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.94"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.26%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "155.36"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.51%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.14"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.94%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.10"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.40%  "

$ws.Range("E35").Value = "  -0.12%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0783"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.81%  "

$ws.Range("E37").Value = "  -0.89%  "

$ws.Range("E38").Value = "  -4.36%  "

$ws.Range("E39").Value = "  -0.78%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.68"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.29%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.32"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +10.81%  "

$ws.Range("E42").Value = "  +0.12%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.78"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.76%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0299"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.50%  "

$ws.Range("E45").Value = "  -6.87%  "

$ws.Range("D46").Value = "2.004.77"
$ws.Range("E46").Value = "  -1.86%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "84.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.08%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.80"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.80%  "

$ws.Range("D49").Value = "2.765.79"
$ws.Range("E49").Value = "  -0.48%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.36"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.19%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "101.91"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.09%  "
